$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Product Code column (A) stays text, not auto-converted to numbers
$ws.Range("A2:A23").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '1159266'
$ws.Range("B2").Value = 'Mens Personalised Classic Silver Steel Spinner Ring'
$ws.Range("D2").Value = 'songsofinkandsteel'
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 'https://www.notonthehighstreet.com/songsofinkandsteel/product/personalised-silver-stainless-steel-wide-spinner-ring'
$ws.Range("G2").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1159266&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 3
$ws.Range("A3").Value = '1218428'
$ws.Range("B3").Value = 'Personalised 18th Birthday Card Wooden Number Gift'
$ws.Range("D3").Value = 'craftheaven'
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 'https://www.notonthehighstreet.com/craftheaven/product/personalised-18th-birthday-card-wooden-number-gift'
$ws.Range("G3").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1218428&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 4
$ws.Range("A4").Value = '1245496'
$ws.Range("B4").Value = 'Solid Perfume Making Kit'
$ws.Range("D4").Value = 'ourhands'
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 'https://www.notonthehighstreet.com/ourhands/product/solid-perfume-making-kit'
$ws.Range("G4").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1245496&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 5
$ws.Range("A5").Value = '1308378'
$ws.Range("B5").Value = 'Women''s White Cotton Nightdress Sleeveless Pink Lizzie'
$ws.Range("D5").Value = 'minilunn'
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 'https://www.notonthehighstreet.com/minilunn/product/women-s-white-cotton-nightdress-sleeveless-pink-lizzie'
$ws.Range("G5").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1308378&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 6
$ws.Range("A6").Value = '1345877'
$ws.Range("B6").Value = 'Ladies Sterling Silver Or Gold Mesh Bracelet'
$ws.Range("D6").Value = 'hurleyburley'
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 'https://www.notonthehighstreet.com/hurleyburley/product/ladies-sterling-silver-mesh-bracelet'
$ws.Range("G6").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1345877&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 7
$ws.Range("A7").Value = '1353406'
$ws.Range("B7").Value = 'F1 Formula One 2025 Calendar Track T Shirt Gift For Him'
$ws.Range("D7").Value = 'nappyhead'
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 'https://www.notonthehighstreet.com/nappyhead/product/formula-1-track-t-shirt-gift-for-him'
$ws.Range("G7").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1353406&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 8
$ws.Range("A8").Value = '1372940'
$ws.Range("B8").Value = '2015 Personalised 10th Tin Wedding Anniversary Poster'
$ws.Range("D8").Value = 'thewordshack'
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 'https://www.notonthehighstreet.com/thewordshack/product/personalised-10th-tin-wedding-anniversary-poster'
$ws.Range("G8").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1372940&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 9
$ws.Range("A9").Value = '1385338'
$ws.Range("B9").Value = 'Mystery Box'
$ws.Range("D9").Value = 'lucysaysido'
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 'https://www.notonthehighstreet.com/lucysaysido/product/mystery-box'
$ws.Range("G9").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1385338&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 10
$ws.Range("A10").Value = '1399796'
$ws.Range("B10").ClearContents() | Out-Null
$ws.Range("D10").ClearContents() | Out-Null
$ws.Range("E10").Value = 2
$ws.Range("F10").ClearContents() | Out-Null
$ws.Range("G10").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1399796&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 11
$ws.Range("A11").Value = '1404697'
$ws.Range("B11").ClearContents() | Out-Null
$ws.Range("D11").ClearContents() | Out-Null
$ws.Range("E11").Value = 2
$ws.Range("F11").ClearContents() | Out-Null
$ws.Range("G11").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1404697&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 12
$ws.Range("A12").Value = '1453979'
$ws.Range("B12").ClearContents() | Out-Null
$ws.Range("D12").ClearContents() | Out-Null
$ws.Range("E12").Value = 2
$ws.Range("F12").ClearContents() | Out-Null
$ws.Range("G12").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1453979&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 13
$ws.Range("A13").Value = '159976'
$ws.Range("B13").ClearContents() | Out-Null
$ws.Range("D13").ClearContents() | Out-Null
$ws.Range("E13").Value = 2
$ws.Range("F13").ClearContents() | Out-Null
$ws.Range("G13").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=159976&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 14
$ws.Range("A14").Value = '458484'
$ws.Range("B14").ClearContents() | Out-Null
$ws.Range("D14").ClearContents() | Out-Null
$ws.Range("E14").Value = 3
$ws.Range("F14").ClearContents() | Out-Null
$ws.Range("G14").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=458484&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 15
$ws.Range("A15").Value = '469358'
$ws.Range("B15").ClearContents() | Out-Null
$ws.Range("D15").ClearContents() | Out-Null
$ws.Range("E15").Value = 3
$ws.Range("F15").ClearContents() | Out-Null
$ws.Range("G15").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=469358&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 16
$ws.Range("A16").Value = '530476'
$ws.Range("B16").ClearContents() | Out-Null
$ws.Range("D16").ClearContents() | Out-Null
$ws.Range("E16").Value = 3
$ws.Range("F16").ClearContents() | Out-Null
$ws.Range("G16").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=530476&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 17
$ws.Range("A17").Value = '608022'
$ws.Range("B17").ClearContents() | Out-Null
$ws.Range("D17").ClearContents() | Out-Null
$ws.Range("E17").Value = 2
$ws.Range("F17").ClearContents() | Out-Null
$ws.Range("G17").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=608022&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 18
$ws.Range("A18").Value = '786481'
$ws.Range("B18").ClearContents() | Out-Null
$ws.Range("D18").ClearContents() | Out-Null
$ws.Range("E18").Value = 2
$ws.Range("F18").ClearContents() | Out-Null
$ws.Range("G18").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=786481&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 19
$ws.Range("A19").Value = '816668'
$ws.Range("B19").ClearContents() | Out-Null
$ws.Range("D19").ClearContents() | Out-Null
$ws.Range("E19").Value = 2
$ws.Range("F19").ClearContents() | Out-Null
$ws.Range("G19").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=816668&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 20
$ws.Range("A20").Value = '876141'
$ws.Range("B20").ClearContents() | Out-Null
$ws.Range("D20").ClearContents() | Out-Null
$ws.Range("E20").Value = 2
$ws.Range("F20").ClearContents() | Out-Null
$ws.Range("G20").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=876141&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 21
$ws.Range("A21").Value = '878687'
$ws.Range("B21").ClearContents() | Out-Null
$ws.Range("D21").ClearContents() | Out-Null
$ws.Range("E21").Value = 2
$ws.Range("F21").ClearContents() | Out-Null
$ws.Range("G21").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=878687&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 22
$ws.Range("A22").Value = '879692'
$ws.Range("B22").ClearContents() | Out-Null
$ws.Range("D22").ClearContents() | Out-Null
$ws.Range("E22").Value = 2
$ws.Range("F22").ClearContents() | Out-Null
$ws.Range("G22").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=879692&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Row 23
$ws.Range("A23").Value = '937471'
$ws.Range("B23").ClearContents() | Out-Null
$ws.Range("D23").ClearContents() | Out-Null
$ws.Range("E23").Value = 2
$ws.Range("F23").ClearContents() | Out-Null
$ws.Range("G23").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=937471&displayFeedbackType=PRODUCT&timeFrame=ALL'

# Remove now-unused rows 24-27 (sheet shrank from 27 to 23 rows)
$ws.Range("A24:G27").Delete() | Out-Null

Write-Output "done"